$wb = $excel.ActiveWorkbook

# ---- Sheet "Login": duplicate row 9 into a newly inserted row 10 ----
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Rows.Item(9).Copy()
$wsLogin.Rows.Item(10).Insert()
$wsLogin.Range("A10").Value = "testT2901_1"
$wsLogin.Range("A10").Select()

# ---- Sheet "Logout": duplicate row 9 into a newly appended row 10 ----
$wsLogout = $wb.Worksheets.Item("Logout")
$wsLogout.Range("A9:F9").Copy($wsLogout.Range("A10:F10"))
$wsLogout.Range("A10").Value = "testT2901_1"
$wsLogout.Range("A10").Select()
